$wb = $excel.ActiveWorkbook

# Insert a brand-new first worksheet named "Sheet1" ahead of "Tom"
# (Worksheets.Add() without args would insert before the active sheet,
# so we pass the first sheet explicitly as the "Before" target.)
$ws = $wb.Worksheets.Add($wb.Worksheets.Item(1))
$ws.Name = "Sheet1"

# Combined client/product/quantity/vendor table - one row per client
# sale, pulled together from the Tom / Bill / Steph sheets.
$data = @(
    @("Clients:", "Products:",       "Quantity:", "Vendors:"),
    @("Tom",      "Leaf Blower",      1,          "Amazon"),
    @("Tom",      "Blender",          3,          "Walmart"),
    @("Tom",      "PS5",              1,          "Gamestop"),
    @("Bill",     "Basketball Hoop",  1,          "Walmart"),
    @("Bill",     "Ipad",             2,          "Amazon"),
    @("Bill",     "Videogame",        4,          "Gamestop"),
    @("Steph",    "Controller",       1,          "Gamestop"),
    @("Steph",    "Sungalsses",       7,          "Amazon"),
    @("Steph",    "Hoodie",           3,          "Walmart")
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $row[$c]
    }
}

# Match the source column-width tweak on the new "Products:" column.
$ws.Columns.Item(2).ColumnWidth = 13.5

# Land the selection/cursor on the cell below the last data row, like the
# author's session, and make sure this new sheet is the active tab.
$ws.Range("D11").Select()
$ws.Activate()
